$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) values
# Price cells are forced to Text format before assignment so that
# numeric-looking strings (e.g. "217.59") are kept as text, matching
# the source data, then formats are cleared again so no stray style
# is left applied to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.200.95"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.659.29"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.74%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5149"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("E8").Value = "  -3.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06444"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.38%  "

$ws.Range("E10").Value = "  -3.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.664.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.296"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.888.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5546"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8059"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.32"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.240.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.93%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.428"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.042"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.756"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1171"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.996"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.06%  "

$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05214"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.253"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.233"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.14%  "

$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.761"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9315"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.373"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.173.23"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +12.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5702"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.48%  "

$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8388"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.680"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.67"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.798.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("E48").Value = "  -3.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.931"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.72%  "

$ws.Range("E51").Value = "  -3.17%  "
